$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.711.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.90%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.794.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.09%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.554'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.35'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.50%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.286'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.91%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0723'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.34%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0932'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.50%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.052.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.08%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.806.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.84%  '

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.05%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.638'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.48%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.738.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.95%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.15%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.61%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0813'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.14%  '

# Row 21
$ws.Range("E21").Value = '  -0.19%  '

# Row 22
$ws.Range("E22").Value = '  +2.67%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.42%  '

# Row 24
$ws.Range("E24").Value = '  +0.00%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.35%  '

# Row 28
$ws.Range("E28").Value = '  -0.14%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0533'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.84%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.16%  '

# Row 32
$ws.Range("E32").Value = '  -1.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.36%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.439.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.59%  '

# Row 36
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.65%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0192'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.50%  '

# Row 38
$ws.Range("E38").Value = '  -0.18%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '84.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.08%  '

# Row 40
$ws.Range("E40").Value = '  -1.93%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.931'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.88%  '

# Row 42
$ws.Range("E42").Value = '  -0.94%  '

# Row 43
$ws.Range("E43").Value = '  +2.97%  '

# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.45%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.11%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0494'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.08%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.949.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '106.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.53%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.98%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0127'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.45%  '
